$wb = $excel.ActiveWorkbook

# --- Sheet3 (Maguskepek): move selection from B5 to C12 (without leaving it the active tab) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate() | Out-Null
$ws3.Range("C12").Select() | Out-Null

# --- Add new sheet "Ellenfelkepek" at the end of the workbook ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Ellenfelkepek"

# Fill in the header + the enemy image paths
$newSheet.Range("A1").Value = "Eleres"
$newSheet.Range("A2").Value = "Images\Ellenseg\ellenseg.png"
$newSheet.Range("A3").Value = "Images\Ellenseg\ellenseg2.png"
$newSheet.Range("A4").Value = "Images\Ellenseg\ellenseg3.png"
$newSheet.Range("A5").Value = "Images\Ellenseg\golem1.png"
$newSheet.Range("A6").Value = "Images\Ellenseg\golem2.png"
$newSheet.Range("A7").Value = "Images\Ellenseg\golem3.png"
$newSheet.Range("A8").Value = "Images\Ellenseg\golem4.png"

# Widen the column so the paths are fully visible
$newSheet.Columns.Item(1).ColumnWidth = 37.333333

# Make the new sheet the active one, with A11 selected
$newSheet.Activate() | Out-Null
$newSheet.Range("A11").Select() | Out-Null
